{"js": "// The document contains two (otherwise identical) exam-ticket header\n// paragraphs of the form \"\u0415\u041a\u0417\u0410\u041c\u0415\u041d\u0410\u0426\u0406\u0419\u041d\u0418\u0419 \u0411\u0406\u041b\u0415\u0422 \u2116 ____N_____ <\u041f\u0406\u0411>\".\n// Only the ticket #7 paragraph (\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447\") is being\n// renamed to \"\u041a\u0443\u043b\u0456\u043d\u0456\u0447\u0435\u0432 \u041c\u0438\u0445\u0430\u0439\u043b\u043e \u0411\u043e\u0440\u0438\u0441\u043e\u0432\u0438\u0447\"; ticket #8's identical text is\n// left untouched. So we must locate that specific paragraph (not just the\n// first text match) before replacing its two name runs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432\") !== -1 && t.indexOf(\"____7_____\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the exam-ticket #7 paragraph containing '\u041c\u0430\u043b\u0430\u0448\u043e\u0432'.\");\n}\n\n// Replace the surname run (\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432\" -> \"\u041a\u0443\u043b\u0456\u043d\u0456\u0447\u0435\u0432\").\nconst surnameHits = target.search(\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432\", { matchCase: true });\nsurnameHits.load(\"items\");\nawait context.sync();\nif (surnameHits.items.length === 0) {\n  throw new Error(\"Surname run '\u041c\u0430\u043b\u0430\u0448\u043e\u0432' not found in target paragraph.\");\n}\nsurnameHits.items[0].insertText(\"\u041a\u0443\u043b\u0456\u043d\u0456\u0447\u0435\u0432\", \"Replace\");\n\n// Replace the given-name/patronymic run (\"  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447\" -> \" \u041c\u0438\u0445\u0430\u0439\u043b\u043e \u0411\u043e\u0440\u0438\u0441\u043e\u0432\u0438\u0447\").\nconst nameHits = target.search(\"  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447\", { matchCase: true });\nnameHits.load(\"items\");\nawait context.sync();\nif (nameHits.items.length === 0) {\n  throw new Error(\"Name run '  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447' not found in target paragraph.\");\n}\nnameHits.items[0].insertText(\" \u041c\u0438\u0445\u0430\u0439\u043b\u043e \u0411\u043e\u0440\u0438\u0441\u043e\u0432\u0438\u0447\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The document contains two (otherwise identical) exam-ticket header\n# paragraphs of the form \"\u0415\u041a\u0417\u0410\u041c\u0415\u041d\u0410\u0426\u0406\u0419\u041d\u0418\u0419 \u0411\u0406\u041b\u0415\u0422 \u2116 ____N_____ <\u041f\u0406\u0411>\".\n# Only the ticket #7 paragraph (\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447\") is being\n# renamed to \"\u041a\u0443\u043b\u0456\u043d\u0456\u0447\u0435\u0432 \u041c\u0438\u0445\u0430\u0439\u043b\u043e \u0411\u043e\u0440\u0438\u0441\u043e\u0432\u0438\u0447\"; ticket #8's identical text is\n# left untouched. So we must locate that specific paragraph (not just the\n# first text match) before replacing its two name runs.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t.Contains(\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432\") -and $t.Contains(\"____7_____\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the exam-ticket #7 paragraph containing '\u041c\u0430\u043b\u0430\u0448\u043e\u0432'.\"\n}\n\n# Replace the surname run (\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432\" -> \"\u041a\u0443\u043b\u0456\u043d\u0456\u0447\u0435\u0432\").\n$surnameRange = $d.Paragraphs($targetIndex).Range\n$surnameRange.Find.Execute(\"\u041c\u0430\u043b\u0430\u0448\u043e\u0432\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u041a\u0443\u043b\u0456\u043d\u0456\u0447\u0435\u0432\", 2)\n\n# Replace the given-name/patronymic run (\"  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447\" -> \" \u041c\u0438\u0445\u0430\u0439\u043b\u043e \u0411\u043e\u0440\u0438\u0441\u043e\u0432\u0438\u0447\").\n$nameRange = $d.Paragraphs($targetIndex).Range\n$nameRange.Find.Execute(\"  \u0413\u0435\u043e\u0440\u0433\u0456\u0439 \u041c\u0438\u043a\u043e\u043b\u0430\u0439\u043e\u0432\u0438\u0447\", $false, $false, $false, $false, $false, $true, 1, $false, \" \u041c\u0438\u0445\u0430\u0439\u043b\u043e \u0411\u043e\u0440\u0438\u0441\u043e\u0432\u0438\u0447\", 2)\n"}
